$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-44 down to 16-45
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record
$ws.Cells.Item(15, 1).Value = 8
$ws.Cells.Item(15, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 45036
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100101
$ws.Cells.Item(15, 8).Value = "Berries"
$ws.Cells.Item(15, 9).Value = 100101001
$ws.Cells.Item(15, 10).Value = "Arándano (blue)"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 10000
$ws.Cells.Item(15, 15).Value = 11000
$ws.Cells.Item(15, 16).Value = 10500
$ws.Cells.Item(15, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(15, 19).Value = 5250
$ws.Cells.Item(15, 20).Value = 2
